# "2nd mission with Vertical layout"
# The shipping country for order PO7375 (row 7) is corrected from "India"
# to "United Kingdom", matching the rest of that shipment's records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "United Kingdom"

# Re-create the row heights that Excel recomputed for the wrapped
# "Shipping Country" column after the edit (rows with the longer
# "United States of America" / "United Kingdom" text wrap to two lines;
# "Germany" rows stay single-line and keep the default height).
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).RowHeight = 38.7
$ws.Rows.Item(3).RowHeight = 25.8
$ws.Rows.Item(4).RowHeight = 25.8
$ws.Rows.Item(5).RowHeight = 25.8
$ws.Rows.Item(6).RowHeight = 25.8
$ws.Rows.Item(7).RowHeight = 25.8
$ws.Rows.Item(11).RowHeight = 38.7

# The author's last selection before saving ended up on I6.
$ws.Range("I6").Select() | Out-Null
